$wb = $excel.ActiveWorkbook

# Update the Users sheet: replace "Brian Miller" with "Jennie Stewart"
$usersSheet = $wb.Worksheets.Item("Users")
$usersSheet.Range("A2").Value = "Jennie Stewart"

# Select cell N12 on the Users sheet and make it the active sheet/tab
$usersSheet.Activate()
$usersSheet.Range("N12").Select()

$wb.Save()
